$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.505.83'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '1.616.76'
$ws.Range("E3").Value = '  -1.48%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.86'
$ws.Range("E5").Value = '  -0.80%  '

$ws.Range("E6").Value = '  -1.79%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.72'
$ws.Range("E8").Value = '  -0.87%  '

$ws.Range("E9").Value = '  +2.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0612'
$ws.Range("E10").Value = '  +0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0886'
$ws.Range("E11").Value = '  -0.27%  '

$ws.Range("D12").Value = '1.844.77'
$ws.Range("E12").Value = '  -1.59%  '

$ws.Range("D13").Value = '1.614.82'
$ws.Range("E13").Value = '  -1.56%  '

$ws.Range("E14").Value = '  -0.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.550'
$ws.Range("E15").Value = '  -2.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.83'
$ws.Range("E16").Value = '  +1.25%  '

$ws.Range("D17").Value = '27.506.64'
$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.43'
$ws.Range("E18").Value = '  +0.90%  '

$ws.Range("D19").Value = '0.0₃0719'
$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("E20").Value = '  -1.62%  '

$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("E22").Value = '  -0.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.09'
$ws.Range("E23").Value = '  +0.80%  '

$ws.Range("E24").Value = '  +7.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.47'
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.81'
$ws.Range("E28").Value = '  -2.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.56'
$ws.Range("E29").Value = '  -0.24%  '

$ws.Range("E30").Value = '  -0.61%  '

$ws.Range("E31").Value = '  -0.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.26'
$ws.Range("E32").Value = '  -1.17%  '

$ws.Range("D33").Value = '1.440.39'
$ws.Range("E33").Value = '  +0.88%  '

$ws.Range("E34").Value = '  -3.62%  '

$ws.Range("E35").Value = '  -3.49%  '

$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("E37").Value = '  +4.94%  '

$ws.Range("E38").Value = '  -2.47%  '

$ws.Range("E39").Value = '  +0.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.859'
$ws.Range("E40").Value = '  -1.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.05'
$ws.Range("E41").Value = '  +6.02%  '

$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("E43").Value = '  -2.69%  '

$ws.Range("E44").Value = '  +0.17%  '

$ws.Range("E45").Value = '  -2.23%  '

$ws.Range("E46").Value = '  -2.77%  '

$ws.Range("D47").Value = '1.755.65'
$ws.Range("E47").Value = '  -1.57%  '

$ws.Range("E48").Value = '  -0.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.42'
$ws.Range("E49").Value = '  +0.13%  '

$ws.Range("E50").Value = '  -1.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0992'
$ws.Range("E51").Value = '  +0.82%  '
